# Mark additional "Test Result" scenarios as updated (TRUE) so the
# summary formula in E1 depends on the actual app settings/status rather
# than staying at the default all-FALSE state.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Result")

# Rows whose STATUS (column C) should flip from FALSE to TRUE.
$rowsToUpdate = @(2, 4, 6, 7, 8, 9, 10, 11)

foreach ($row in $rowsToUpdate) {
    $ws.Cells.Item($row, 3).Value = $true
}
